# Applies weekly-rotation update to the "Ajo" (garlic) price sheet.
# Columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are updated per row
# according to the source data diff. Row 20 and header row 1 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(D, J, K, L, M, P)
$data = @{
    2  = @(44160,  360, 10000, 11000, 10500, 1050)
    3  = @(44727,  400, 18000, 19000, 18500, 1850)
    4  = @(44890,  400, 16000, 17000, 16500, 1650)
    5  = @(44679,  200, 19000, 20000, 19500, 1950)
    6  = @(44218,  320, 10000, 11000, 10500, 1050)
    7  = @(44291,  200, 13000, 14000, 13500, 1350)
    8  = @(44263,  300, 15000, 16000, 15500, 1550)
    9  = @(44777,  200, 24000, 25000, 24500, 2450)
    10 = @(44580,  200, 18000, 20000, 19000, 1900)
    11 = @(44358,  300, 14000, 15000, 14500, 1450)
    12 = @(44330,  300, 13000, 14000, 13500, 1350)
    13 = @(44406,  400, 14000, 15000, 14500, 1450)
    14 = @(44204,  400, 10000, 11000, 10500, 1050)
    15 = @(44860,  400, 14000, 15000, 14500, 1450)
    16 = @(44441,  300, 15000, 16000, 15500, 1550)
    17 = @(44547,  300, 19000, 20000, 19500, 1950)
    18 = @(44882,  400, 15000, 16000, 15550, 1555)
    19 = @(44694,  400, 16000, 17000, 16500, 1650)
    21 = @(44714,  400, 19000, 20000, 19500, 1950)
    22 = @(44428,  300, 15000, 16000, 15500, 1550)
    23 = @(44847,  400, 16000, 17000, 16500, 1650)
    24 = @(44377,  650, 14000, 15000, 14538, 1454)
    25 = @(44893, 1400, 15000, 16000, 15571, 1557)
    26 = @(44524,  200, 20000, 21000, 20500, 2050)
    27 = @(44644,  300, 20000, 21000, 20500, 2050)
    28 = @(44460,  300, 15000, 16000, 15500, 1550)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
    $ws.Range("K$row").Value = $vals[2]
    $ws.Range("L$row").Value = $vals[3]
    $ws.Range("M$row").Value = $vals[4]
    $ws.Range("P$row").Value = $vals[5]
}
